# add mysql connect info
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 5: "Pwd" -> "SqlIP" -------------------------------------------------
$ws.Range("A5").Value = "SqlIP"
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)

# --- Clone the row-4 formatting onto the three new rows, then fill values --
$ws.Range("A4:J4").Copy()
$ws.Range("A6:J6").PasteSpecial(-4122)
$ws.Range("A4:J4").Copy()
$ws.Range("A7:J7").PasteSpecial(-4122)
$ws.Range("A4:J4").Copy()
$ws.Range("A8:J8").PasteSpecial(-4122)

# --- Row 6: SqlPort / int ----------------------------------------------------
$ws.Range("A6").Value = "SqlPort"
$ws.Range("B6").Value = "int"
$ws.Range("C6").Value = $true
$ws.Range("D6").Value = $false
$ws.Range("E6").Value = $false
$ws.Range("F6").Value = $true
$ws.Range("G6").Value = 0
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = "Friend"

# --- Row 7: SqlUser / string --------------------------------------------------
$ws.Range("A7").Value = "SqlUser"
$ws.Range("B7").Value = "string"
$ws.Range("C7").Value = $true
$ws.Range("D7").Value = $false
$ws.Range("E7").Value = $false
$ws.Range("F7").Value = $true
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = "Friend"
$ws.Range("A7").ClearFormats()
$ws.Range("A7").Value = "SqlUser"

# --- Row 8: SqlPwd / string ----------------------------------------------------
$ws.Range("A8").Value = "SqlPwd"
$ws.Range("B8").Value = "string"
$ws.Range("C8").Value = $true
$ws.Range("D8").Value = $false
$ws.Range("E8").Value = $false
$ws.Range("F8").Value = $true
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = "Friend"
$ws.Range("A8").ClearFormats()
$ws.Range("A8").Value = "SqlPwd"

# --- Selection / active cell, as left by the author after editing ----------
$ws.Range("C15").Select()
